# 6 lesson Done HW
# Append 33 new test-result rows (rows 35..67) to Sheet1, mirroring the
# existing data layout (A: seq#, B: test name, C/D: pass/fail marker, E: timestamp).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Template cells whose number/font formatting we reuse for the new rows,
# so no new (duplicate) cell styles get created in the workbook.
$successTemplate = $ws.Cells.Item(2, 3)   # green "success" style
$failureTemplate = $ws.Cells.Item(3, 4)   # red/bold "failure" style
$dateTemplate    = $ws.Cells.Item(2, 5)   # date/time style

$names = @{
  5  = "TestRebelStar::test_domain"
  6  = "TestRebelStar::test_method[True-True-True]"
  7  = "TestRebelStar::test_method[True-True-False]"
  8  = "TestRebelStar::test_method[True-False-True]"
  9  = "TestRebelStar::test_method[False-True-True]"
  10 = "TestRebelStar::test_method[True-False-False]"
  11 = "TestRebelStar::test_method[False-True-False]"
  12 = "TestRebelStar::test_method[False-False-True]"
  13 = "TestRebelStar::test_method[False-False-False]"
  14 = "TestRebelStar::test_crush_method"
  15 = "TestRebelStar::test_crush_method_1"
  16 = "TestRebelStar::test_crush_method_2"
}

# Each entry: row, A-value(seq#), name-key, outcome-column(C=success,D=failure), outcome-value, E-timestamp
$data = @(
  @(35,34,5,"C",1,44087.81463707585),
  @(36,35,5,"D",1,44087.81464879269),
  @(37,36,5,"C",1,44087.81464881967),
  @(38,37,6,"C",1,44087.81464897447),
  @(39,38,6,"C",1,44087.81465225227),
  @(40,39,6,"C",1,44087.81465226878),
  @(41,40,7,"C",1,44087.81465229909),
  @(42,41,7,"C",1,44087.81465525642),
  @(43,42,7,"C",1,44087.81465527238),
  @(44,43,8,"C",1,44087.81465529553),
  @(45,44,8,"C",1,44087.81465877905),
  @(46,45,8,"C",1,44087.81465879554),
  @(47,46,9,"C",1,44087.81465881842),
  @(48,47,9,"C",1,44087.81466267136),
  @(49,48,9,"C",1,44087.81466269307),
  @(50,49,10,"C",1,44087.81466271514),
  @(51,50,10,"C",1,44087.81466556362),
  @(52,51,10,"C",1,44087.81466557922),
  @(53,52,11,"C",1,44087.81466560782),
  @(54,53,11,"C",1,44087.81466854236),
  @(55,54,11,"C",1,44087.81466855935),
  @(56,55,12,"C",1,44087.81466858118),
  @(57,56,12,"C",1,44087.81467160763),
  @(58,57,12,"C",1,44087.8146716242),
  @(59,58,13,"C",1,44087.81467164343),
  @(60,59,13,"C",1,44087.81467425697),
  @(61,60,13,"C",1,44087.81467427583),
  @(62,61,14,"D",1,44087.81467433368),
  @(63,62,14,"C",1,44087.81467437302),
  @(64,63,15,"D",1,44087.81467439864),
  @(65,64,15,"C",1,44087.81467441384),
  @(66,65,16,"D",1,44087.81467443667),
  @(67,66,16,"C",1,44087.81467445286)
)

foreach ($row in $data) {
    $r       = $row[0]
    $seqNum  = $row[1]
    $nameKey = $row[2]
    $col     = $row[3]
    $val     = $row[4]
    $ts      = $row[5]

    $ws.Cells.Item($r, 1).Value = $seqNum
    $ws.Cells.Item($r, 2).Value = $names[$nameKey]

    if ($col -eq "D") {
        $outcome = $ws.Cells.Item($r, 4)
        $outcome.Value = $val
        $failureTemplate.Copy()
        $outcome.PasteSpecial($xlPasteFormats)
    } else {
        $outcome = $ws.Cells.Item($r, 3)
        $outcome.Value = $val
        $successTemplate.Copy()
        $outcome.PasteSpecial($xlPasteFormats)
    }

    $tsCell = $ws.Cells.Item($r, 5)
    $tsCell.Value = $ts
    $dateTemplate.Copy()
    $tsCell.PasteSpecial($xlPasteFormats)
}

$excel.CutCopyMode = $false
